# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold + border) from the "Weekly Quantity" header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-serial number format from the "Weekly Quantity" date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A24").PasteSpecial(-4122)

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows

$wsForecast.Range("A2").Value = 45361.99999999999
$wsForecast.Range("B2").Value = 32
$wsForecast.Range("C2").Value = -3.647613261870044
$wsForecast.Range("D2").Value = 66.54810955440122
$wsForecast.Range("A3").Value = 45368.99999999999
$wsForecast.Range("B3").Value = 31
$wsForecast.Range("C3").Value = -2.55162465506946
$wsForecast.Range("D3").Value = 65.83365000431347
$wsForecast.Range("A4").Value = 45389.99999999999
$wsForecast.Range("B4").Value = 29
$wsForecast.Range("C4").Value = -4.409602393098049
$wsForecast.Range("D4").Value = 61.76934435490998
$wsForecast.Range("A5").Value = 45410.99999999999
$wsForecast.Range("B5").Value = 26
$wsForecast.Range("C5").Value = -4.680457288908126
$wsForecast.Range("D5").Value = 59.2262730941695
$wsForecast.Range("A6").Value = 45466.99999999999
$wsForecast.Range("B6").Value = 20
$wsForecast.Range("C6").Value = -13.59076290354936
$wsForecast.Range("D6").Value = 53.85647201316866
$wsForecast.Range("A7").Value = 45480.99999999999
$wsForecast.Range("B7").Value = 18
$wsForecast.Range("C7").Value = -17.01972985064119
$wsForecast.Range("D7").Value = 52.26429499121329
$wsForecast.Range("A8").Value = 45487.99999999999
$wsForecast.Range("B8").Value = 18
$wsForecast.Range("C8").Value = -16.26812035281074
$wsForecast.Range("D8").Value = 52.11857018321028
$wsForecast.Range("A9").Value = 45494.99999999999
$wsForecast.Range("B9").Value = 17
$wsForecast.Range("C9").Value = -17.16419910156114
$wsForecast.Range("D9").Value = 49.43437717821867
$wsForecast.Range("A10").Value = 45501.99999999999
$wsForecast.Range("B10").Value = 16
$wsForecast.Range("C10").Value = -17.90151737256482
$wsForecast.Range("D10").Value = 48.9517819764591
$wsForecast.Range("A11").Value = 45529.99999999999
$wsForecast.Range("B11").Value = 13
$wsForecast.Range("C11").Value = -20.35942035732005
$wsForecast.Range("D11").Value = 48.03478196263231
$wsForecast.Range("A12").Value = 45543.99999999999
$wsForecast.Range("B12").Value = 11
$wsForecast.Range("C12").Value = -22.63266545524574
$wsForecast.Range("D12").Value = 46.05938507686926
$wsForecast.Range("A13").Value = 45571.99999999999
$wsForecast.Range("B13").Value = 8
$wsForecast.Range("C13").Value = -24.53698153287244
$wsForecast.Range("D13").Value = 39.23206014624527
$wsForecast.Range("A14").Value = 45578.99999999999
$wsForecast.Range("B14").Value = 7
$wsForecast.Range("C14").Value = -26.14506875083705
$wsForecast.Range("D14").Value = 40.84174249848795
$wsForecast.Range("A15").Value = 45585.99999999999
$wsForecast.Range("B15").Value = 6
$wsForecast.Range("C15").Value = -27.10269197962714
$wsForecast.Range("D15").Value = 40.98937165866118
$wsForecast.Range("A16").Value = 45613.99999999999
$wsForecast.Range("B16").Value = 3
$wsForecast.Range("C16").Value = -30.15188755862124
$wsForecast.Range("D16").Value = 37.95374650962286
$wsForecast.Range("A17").Value = 45620.99999999999
$wsForecast.Range("B17").Value = 2
$wsForecast.Range("C17").Value = -28.26265200670499
$wsForecast.Range("D17").Value = 37.33349580052481
$wsForecast.Range("A18").Value = 45627.99999999999
$wsForecast.Range("B18").Value = 2
$wsForecast.Range("C18").Value = -32.72800476818869
$wsForecast.Range("D18").Value = 35.48234424124013
$wsForecast.Range("A19").Value = 45634.99999999999
$wsForecast.Range("B19").Value = 1
$wsForecast.Range("C19").Value = -29.27618046293266
$wsForecast.Range("D19").Value = 37.86092498992465
$wsForecast.Range("A20").Value = 45641.99999999999
$wsForecast.Range("B20").Value = 0
$wsForecast.Range("C20").Value = -32.47663317836107
$wsForecast.Range("D20").Value = 34.01282007526821
$wsForecast.Range("A21").Value = 45648.99999999999
$wsForecast.Range("B21").Value = 0
$wsForecast.Range("C21").Value = -35.53825343871706
$wsForecast.Range("D21").Value = 32.71416625148449
$wsForecast.Range("A22").Value = 45655.99999999999
$wsForecast.Range("B22").Value = 0
$wsForecast.Range("C22").Value = -36.57664830248579
$wsForecast.Range("D22").Value = 34.72689612698971
$wsForecast.Range("A23").Value = 45662.99999999999
$wsForecast.Range("B23").Value = 0
$wsForecast.Range("C23").Value = -33.15363903374651
$wsForecast.Range("D23").Value = 29.380806677137
$wsForecast.Range("A24").Value = 45669.99999999999
$wsForecast.Range("B24").Value = 0
$wsForecast.Range("C24").Value = -34.9826439055968
$wsForecast.Range("D24").Value = 32.22509418173549

$wsForecast.Range("A1").Select()

